$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet1) ---
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B): 5.0.0 -> 6.0.0
$ws.Cells.Item(3, 2).Value = "6.0.0"

# Update Date value (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Row 9 is "Publisher" with an empty value -> set it to "Alvearie Team"
$ws.Cells.Item(9, 2).Value = "Alvearie Team"

# Rows 10 and 11 were a duplicated "Contact" / "No display for ContactDetail"
# pair. Repurpose row 10 into the new "Jurisdiction" property and row 11
# into the "Description" property (previously row 12).
$ws.Cells.Item(10, 1).Value = "Jurisdiction"
$ws.Cells.Item(10, 2).Value = "United States of America"
$ws.Cells.Item(11, 1).Value = "Description"
$ws.Cells.Item(11, 2).Value = "Encrypted state"

# The original "Description" row (now redundant) is removed, netting the
# sheet down from 21 to 20 rows.
$ws.Rows.Item(12).Delete()

# --- Sheet "Elements" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Elements")

# Root Extension element row (row 2): Short/Definition now reflect the
# StructureDefinition's own Title/Description instead of the generic
# "Extension" / "An Extension" placeholders.
$ws2.Cells.Item(2, 11).Value = "Encrypted State"
$ws2.Cells.Item(2, 12).Value = "Encrypted state"
